$wb = $excel.ActiveWorkbook

# --- Part 1: update B102 ("remn_amt" for the last date row) on the first 6 sheets ---
$newB102 = @(3193, 1163, 1296, 1893, 758, 1510)
for ($i = 1; $i -le 6; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Range("B102").Value = $newB102[$i - 1]
}

# --- Part 2: add a new sheet "대영포장" (7th sheet) with a date/remn_amt table ---
$ws1 = $wb.Worksheets.Item(1)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "대영포장"

# Clone header + row formatting (header style + date number format) from sheet1's
# first 101 rows (header row + 100 data rows), then overwrite with this sheet's own values.
$ws1.Range("A1:B101").Copy($newSheet.Range("A1"))

$dates = @(45813,45817,45818,45819,45820,45821,45824,45825,45826,45827,45828,45831,45832,45833,45834,45835,45838,45839,45840,45841,45842,45845,45846,45847,45848,45849,45852,45853,45854,45855,45856,45859,45860,45861,45862,45863,45866,45867,45868,45869,45870,45873,45874,45875,45876,45877,45880,45881,45882,45883,45887,45888,45889,45890,45891,45894,45895,45896,45897,45898,45901,45902,45903,45904,45905,45908,45909,45910,45911,45912,45915,45916,45917,45918,45919,45922,45923,45924,45925,45926,45929,45930,45931,45932,45940,45943,45944,45945,45946,45947,45950,45951,45952,45953,45954,45957,45958,45959,45960,45961)
$vals  = @(3138,3253,3176,3486,3625,3708,3706,3766,3747,3753,3807,3559,3548,3522,3482,3509,3512,3539,3584,3739,3601,3627,3651,3887,3886,3948,3863,3876,3844,3841,3916,3919,3840,3802,3722,3719,3715,3823,3859,3984,3828,3878,3978,4085,4263,4097,4059,3838,3904,3882,3937,3943,3900,3923,3943,3945,3911,3842,3850,3845,3842,3781,3789,3876,3719,3679,3347,3329,3409,3374,3388,3363,3341,3372,3310,3285,3288,3232,3243,3173,3170,3145,3149,3190,3154,3122,3106,3187,3151,3106,3103,3111,3117,3066,3032,3035,3019,2964,2911,2911)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $newSheet.Cells.Item($row, 1).Value = $dates[$i]
    $newSheet.Cells.Item($row, 2).Value = $vals[$i]
}

Write-Output "done"
